# Opp Final - 31 Mar 2025
# Correct the spelling of "Rachell Schaller" -> "Rachel Schaller" on the
# Users sheet, then leave the Users sheet active (with cell F9 selected),
# matching the editor's final on-screen state after making the change.

$wb = $excel.ActiveWorkbook

$usersSheet = $wb.Worksheets.Item("Users")
$usersSheet.Range("B2").Value = "Rachel Schaller"

$usersSheet.Activate()
$usersSheet.Range("F9").Select()
